# Updated POL model - 2025-08-19 17:53
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the first sheet from "Sheet1" to "misc."
$ws.Name = "misc."

# Add two new header columns (K4, L4) matching the existing header style (same as J4)
$ws.Range("K4").Value = "other_indexes"
$ws.Range("K4").Style = $ws.Range("J4").Style
$ws.Range("L4").Value = "commodity"
$ws.Range("L4").Style = $ws.Range("J4").Style

# New row 11: flo_emis / gas / *ccs,*ccs-rf / 0.95 / co2 / co2captured
$ws.Range("B11").Value = "flo_emis"
$ws.Range("D11").Value = "gas"
$ws.Range("L11").Value = "co2captured"
$ws.Range("K11").Value = "co2"
$ws.Range("E11").Value = "*ccs,*ccs-rf"
$ws.Range("H11").Value = 0.95

# New row 12: flo_emis / coal,oil / *ccs,*ccs-rf / 0.85 / co2 / co2captured
$ws.Range("B12").Value = "flo_emis"
$ws.Range("D12").Value = "coal,oil"
$ws.Range("E12").Value = "*ccs,*ccs-rf"
$ws.Range("H12").Value = 0.85
$ws.Range("K12").Value = "co2"
$ws.Range("L12").Value = "co2captured"

# Cosmetic: column widths (best-fit approximations) for the touched columns
$ws.Columns.Item(5).ColumnWidth = 9.53125
$ws.Columns.Item(11).ColumnWidth = 11.28125

# Move the active selection to D13, matching the post-edit cursor position
$ws.Range("D13").Select() | Out-Null

Write-Output "done"
